# Updated cryptos list on Sat Feb 24 17:51:44 UTC 2024 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns for every coin row, and
# re-sorts a handful of rows whose relative ranking changed (Stellar/Celestia,
# Stacks/ARBITRUM, WEMIXToken/TheGraph swap places).
#
# Price cells that look like plain numbers ("1.00", "7.82", "0.0444", ...)
# are forced to Text format first so Excel doesn't silently coerce them to
# numeric values (losing the trailing zero / leading-zero formatting used
# throughout this sheet). Cells whose new price contains two dots
# (e.g. "51.581.09") already fail numeric parsing, so no extra handling is
# required for those.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.581.09'
$ws.Range("E2").Value = '  +0.95%  '
$ws.Range("D3").Value = '2.994.66'
$ws.Range("E3").Value = '  +1.77%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '381.95'
$ws.Range("E5").Value = '  +2.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.65'
$ws.Range("E6").Value = '  +2.35%  '
$ws.Range("E7").Value = '  +2.42%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.596'
$ws.Range("E9").Value = '  +2.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.74'
$ws.Range("E10").Value = '  +1.24%  '
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("E12").Value = '  +1.63%  '
$ws.Range("D13").Value = '3.463.76'
$ws.Range("E13").Value = '  +1.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.84'
$ws.Range("E14").Value = '  +3.96%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '18.47'
$ws.Range("E15").Value = '  +2.44%  '
$ws.Range("D16").Value = '2.998.65'
$ws.Range("E16").Value = '  +1.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '11.24'
$ws.Range("E17").Value = '  +1.94%  '
$ws.Range("E18").Value = '  +2.37%  '
$ws.Range("D19").Value = '51.586.50'
$ws.Range("E19").Value = '  +1.01%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.14'
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.62'
$ws.Range("E21").Value = '  +1.57%  '
$ws.Range("E22").Value = '  +0.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.38'
$ws.Range("E23").Value = '  +2.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '268.23'
$ws.Range("E24").Value = '  +1.19%  '
$ws.Range("E25").Value = '  +3.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.82'
$ws.Range("E26").Value = '  -4.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.47'
$ws.Range("E27").Value = '  -1.18%  '
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '26.16'
$ws.Range("E29").Value = '  +1.97%  '
$ws.Range("E30").Value = '  +1.57%  '
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.37'
$ws.Range("E32").Value = '  +3.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.90'
$ws.Range("E33").Value = '  +4.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '51.52'
$ws.Range("E34").Value = '  +1.63%  '
$ws.Range("E35").Value = '  +1.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0444'
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("E38").Value = '  +2.74%  '
$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.117'
$ws.Range("E39").Value = '  +2.03%  '
$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '16.80'
$ws.Range("E40").Value = '  +2.86%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.57'
$ws.Range("E41").Value = '  +3.14%  '
$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.85'
$ws.Range("E42").Value = '  +3.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '125.04'
$ws.Range("E43").Value = '  +4.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.76'
$ws.Range("E45").Value = '  +1.71%  '
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.03'
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("B47").Value = 'TheGraph'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.272'
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("E48").Value = '  +2.97%  '
$ws.Range("D49").Value = '2.042.05'
$ws.Range("E49").Value = '  +2.28%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0334'
$ws.Range("E50").Value = '  +2.73%  '
$ws.Range("E51").Value = '  +16.06%  '
